# ToDo files from chapter 12 to 20
# Change all "X" status markers in column G (chapters 16-20) to "ToDo",
# and fix the two "esthetic differences" name-mismatch rows so the
# R-script/folder name matches the book convention (swap book/github names).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("List")

# Mark remaining "X" entries in the status column as "ToDo"
$rng = $ws.Range("G1:G1048576")
$rng.Replace("X", "ToDo") | Out-Null

# Row 43: MVAportfolIBMFord / MVAportfol_IBM_Ford naming note
$ws.Range("C43").Value = "MVAportfol_IBM_Ford"
$ws.Range("K43").Value = "esthetic differences, name in book: MVAportfolIBMFord"

# Row 57: MVAbancrupcydis / MVAbankruptcydis naming note
$ws.Range("C57").Value = "MVAbankruptcydis"
$ws.Range("K57").Value = "esthetic differences, name in book: MVAbancrupcydis"

# Restore the view: scrolled to top (below the frozen header), selection on A25
$ws.Range("A25").Select()
